$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 53, shifting old rows 53-55 down to 54-56
$ws.Rows.Item(53).Insert()

# Fill in the new row 53 with data (copy common fields, new price/volume/date)
$ws.Cells.Item(53, 1).Value = 6
$ws.Cells.Item(53, 2).Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Cells.Item(53, 3).Value = "Metropolitana"
$ws.Cells.Item(53, 4).Value = 44461
$ws.Cells.Item(53, 5).Value = 13
$ws.Cells.Item(53, 6).Value = "Fruta"
$ws.Cells.Item(53, 7).Value = 100108
$ws.Cells.Item(53, 8).Value = "Tropicales y subtropicales"
$ws.Cells.Item(53, 9).Value = 100108007
$ws.Cells.Item(53, 10).Value = "Coco"
$ws.Cells.Item(53, 11).Value = "Sin especificar"
$ws.Cells.Item(53, 12).Value = "Primera"
$ws.Cells.Item(53, 13).Value = 60
$ws.Cells.Item(53, 14).Value = 19000
$ws.Cells.Item(53, 15).Value = 20000
$ws.Cells.Item(53, 16).Value = 19500
$ws.Cells.Item(53, 17).Value = "$/malla 20 unidades"
$ws.Cells.Item(53, 18).Value = "Perú"
$ws.Cells.Item(53, 19).Value = 975
$ws.Cells.Item(53, 20).Value = 20
